$d = $word.ActiveDocument

$table = $d.Tables.Item(1)

# Row 6 (1-based) is the first empty data row -> 10.02.2023 / 0,5h / Sprintti palaveri
$table.Cell(6, 1).Range.Text = "10.02.2023"
$table.Cell(6, 2).Range.Text = "0,5h"
$table.Cell(6, 3).Range.Text = "Sprintti palaveri"

# Row 7 (1-based) is the second empty data row -> 14.02.2023 / 0,5h / Käyttäjänluomis-ikkunan hiomista ja pohdintaa
$table.Cell(7, 1).Range.Text = "14.02.2023"
$table.Cell(7, 2).Range.Text = "0,5h"
$table.Cell(7, 3).Range.Text = "Käyttäjänluomis-ikkunan hiomista ja pohdintaa"
